$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings need an explicit
# Text number format, otherwise Excel auto-converts them to numbers (e.g.
# "0.850" -> 0.85, losing the trailing zero / exact text representation).

$ws.Range("D2").Value = "60.708.72"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.590.95"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.77"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.57"
$ws.Range("E6").Value = "  -2.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  -4.45%  "
$ws.Range("D9").Value = "2.596.95"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.54"
$ws.Range("E10").Value = "  +6.54%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "3.040.26"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "60.604.38"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.73"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "2.585.56"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.63"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.48"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.13"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "2.701.32"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "0.0₃0850"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.42"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.43"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "153.24"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.73"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.01"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.854"
$ws.Range("E38").Value = "  +6.85%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.850"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.19"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.76"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "296.32"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.623"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0995"
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0560"
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.83"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.89"
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0234"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").Value = "  +0.48%  "
